$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B29").Value = "Salto De Agua"
$ws.Range("B30").Value = "San Cristóbal De Las Casas"
$ws.Range("B52").Value = "Coyame Del Sotol"
$ws.Range("B61").Value = "Guadalupe Y Calvo"
$ws.Range("B64").Value = "Hidalgo Del Parral"
$ws.Range("B78").Value = "San Francisco De Conchos"
$ws.Range("B79").Value = "San Francisco Del Oro"
$ws.Range("B84").Value = "Valle De Zaragoza"
$ws.Range("A102").Value = "Ciudad De México"
$ws.Range("B106").Value = "Cuajimalpa De Morelos"
$ws.Range("B118").Value = "Coneto De Comonfort"
$ws.Range("B132").Value = "Nombre De Dios"
$ws.Range("B136").Value = "Pánuco De Coronado"
$ws.Range("B142").Value = "San Juan De Guadalupe"
$ws.Range("B143").Value = "San Juan Del Río"
$ws.Range("B144").Value = "San Luis Del Cordero"
$ws.Range("B145").Value = "San Pedro Del Gallo"
$ws.Range("A155").Value = "Estado De México"
$ws.Range("B155").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B159").Value = "Atizapán De Zaragoza"
$ws.Range("B165").Value = "Coacalco De Berriozábal"
$ws.Range("B166").Value = "Ecatepec De Morelos"
$ws.Range("B169").Value = "Naucalpan De Juárez"
$ws.Range("B171").Value = "San Martín De Las Pirámides"
$ws.Range("B176").Value = "Tenango Del Valle"
$ws.Range("B179").Value = "Tlalnepantla De Baz"
$ws.Range("B183").Value = "Villa De Allende"
$ws.Range("A188").Value = "Guanajuato"
$ws.Range("B191").Value = "Apaseo El Alto"
$ws.Range("B192").Value = "Apaseo El Grande"
$ws.Range("B198").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B201").Value = "Jaral Del Progreso"
$ws.Range("B207").Value = "Purísima Del Rincón"
$ws.Range("B212").Value = "San Francisco Del Rincón"
$ws.Range("B214").Value = "San Luis De La Paz"
$ws.Range("B215").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B216").Value = "Silao De La Victoria"
$ws.Range("B219").Value = "Valle De Santiago"
$ws.Range("B221").Value = "Acapulco De Juárez"
$ws.Range("B224").Value = "Ajuchitlán Del Progreso"
$ws.Range("B225").Value = "Alcozauca De Guerero"
$ws.Range("B227").Value = "Atenango Del Río"
$ws.Range("B228").Value = "Atoyac De Álvarez"
$ws.Range("B229").Value = "Ayutla De Los Libres"
$ws.Range("B232").Value = "Chilapa De Álvarez"
$ws.Range("B233").Value = "Chilpancingo De Los Bravo"
$ws.Range("B234").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B238").Value = "Coyuca De Benítez"
$ws.Range("B239").Value = "Coyuca De Catalán"
$ws.Range("B242").Value = "Cutzamala De Pinzón"
$ws.Range("B244").Value = "Iguala De La Independencia"
$ws.Range("B254").Value = "Taxco De Alarcón"
$ws.Range("B256").Value = "Técpan De Galeana"
$ws.Range("B261").Value = "Tlapa De Comonfort"
$ws.Range("B280").Value = "Mineral Del Chico"
$ws.Range("B281").Value = "Mixquiahuala De Juárez"
$ws.Range("B282").Value = "Molango De Escamilla"
$ws.Range("B283").Value = "Pachuca De Soto"
$ws.Range("B284").Value = "Progreso De Obregón"
$ws.Range("B285").Value = "Santiago Tulantepec De Lugo Guerero"
$ws.Range("B288").Value = "Tenango De Doria"
$ws.Range("B290").Value = "Tepehuacán De Guerero"
$ws.Range("B291").Value = "Tezontepec De Aldama"
$ws.Range("B293").Value = "Tula De Allende"
$ws.Range("B294").Value = "Tulancingo De Bravo"
$ws.Range("B299").Value = "Ahualulco De Mercado"
$ws.Range("B302").Value = "Atotonilco El Alto"
$ws.Range("B310").Value = "Encarnación De Díaz"
$ws.Range("B317").Value = "Lagos De Moreno"
$ws.Range("B322").Value = "San Juan De Los Lagos"
$ws.Range("B323").Value = "San Juanito De Escobedo"
$ws.Range("B327").Value = "San Miguel El Alto"
$ws.Range("B329").Value = "Tamazula De Gordiano"
$ws.Range("B331").Value = "Teocuitatlán De Corona"
$ws.Range("B332").Value = "Tepatitlán De Morelos"
$ws.Range("B338").Value = "Unión De San Antonio"
$ws.Range("B340").Value = "Yahualica De González Gallo"
$ws.Range("B342").Value = "Zapotitlán De Vadillo"
$ws.Range("B402").Value = "Ixtlán Del Río"
$ws.Range("B416").Value = "Montemorelos"
$ws.Range("B418").Value = "San Nicolás De Los Garza"
$ws.Range("B420").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B421").Value = "Ayoquezco De Aldama"
$ws.Range("B422").Value = "Constancia Del Rosario"
$ws.Range("B424").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B425").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B426").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B429").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B430").Value = "Oaxaca De Juárez"
$ws.Range("B435").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B462").Value = "Teotitlán De Flores Magón"
$ws.Range("B463").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B477").Value = "Palmar De Bravo"
$ws.Range("B482").Value = "San Salvador El Verde"
$ws.Range("B484").Value = "Tecali De Herrera"
$ws.Range("B485").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B489").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B496").Value = "Amealco De Bonfil"
$ws.Range("B498").Value = "Cadereyta De Montes"
$ws.Range("B501").Value = "Jalpan De Serra"
$ws.Range("B504").Value = "San Juan Del Río"
$ws.Range("B515").Value = "Cerro De San Pedro"
$ws.Range("B517").Value = "Ciudad Del Maíz"
$ws.Range("B521").Value = "Mexquitic De Carmona"
$ws.Range("B526").Value = "Santa María Del Río"
$ws.Range("B533").Value = "Villa De Ramos"
$ws.Range("B577").Value = "Soto La Marina"
$ws.Range("B587").Value = "Papalotla De Xicohténcatl"
$ws.Range("B589").Value = "Tepetitla De Lardizábal"
$ws.Range("B600").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B604").Value = "Boca Del Río"
$ws.Range("B605").Value = "Castillo De Teayo"
$ws.Range("B612").Value = "Cosamaloapan De Carpio"
$ws.Range("B622").Value = "Ixhuatlán De Madero"
$ws.Range("B625").Value = "Juchique De Ferrer"
$ws.Range("B627").Value = "Lerdo De Tejada"
$ws.Range("B628").Value = "Martínez De La Torre"
$ws.Range("B639").Value = "Soledad De Doblado"
$ws.Range("B655").Value = "Cañitas De Felipe Pescador"
$ws.Range("B657").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B671").Value = "Nochistlán De Mejía"
$ws.Range("B672").Value = "Noria De Ángeles"
$ws.Range("B680").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B683").Value = "Villa De Cos"

# --- Floating point precision fixes ---
$ws.Range("D203").Value = 0.009171507184347297
$ws.Range("D330").Value = 0.009171507184347297
$ws.Range("D436").Value = 0.009171507184347297

# --- Remove trailing footer/metadata rows 691-695 (dimension shrinks to A1:D689) ---
$ws.Rows("691:695").Delete()
